# Add a new "blog post" entry (md + image file names) to the "news" sheet
# and make that sheet the active one (matching the author's last on-screen
# state when the workbook was saved).

$wb = $excel.ActiveWorkbook

$newsSheet = $wb.Worksheets.Item("news")

# Fill in the new columns for the existing news row.
$newsSheet.Range("F2").Value = "hello_world.md"
$newsSheet.Range("G2").Value = "hello_word.png"

# Make "news" the active sheet/tab, and leave the selection where the
# author left it after typing the new values.
$newsSheet.Activate()
$newsSheet.Range("G3").Select()
